$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 532, shifting existing rows 532:606 down to 533:607
$ws.Rows(532).Insert()

# Populate the new row 532 with the new observation (constants copied from the
# surrounding rows of this same data block, only date/volume/price columns differ)
$ws.Range("A532").Value = 3
$ws.Range("B532").Value = "Femacal de La Calera"
$ws.Range("C532").Value = "Coquimbo"
$ws.Range("D532").Value = 45127
$ws.Range("E532").Value = 5
$ws.Range("F532").Value = 100112040
$ws.Range("G532").Value = "Cilantro"
$ws.Range("H532").Value = "Sin especificar"
$ws.Range("I532").Value = "Primera"
$ws.Range("J532").Value = 260
$ws.Range("K532").Value = 3800
$ws.Range("L532").Value = 4000
$ws.Range("M532").Value = 3885
$ws.Range("N532").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O532").Value = "Provincia de Quillota"
$ws.Range("P532").Value = 1295
$ws.Range("Q532").Value = 3
$ws.Range("R532").Value = "Hortaliza"
